$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying the existing "2022-Q2"
#    sheet (so it inherits the same header/row styling) and placing
#    the copy immediately before "2022-Q2".
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q2")
$src.Copy($src)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Remove the extra template rows (the source sheet had 6 data rows,
# the new one only needs 3).
$newSheet.Rows("5:7").Delete()

# ------------------------------------------------------------------
# 2. Populate "2022-Q3" with its fund holdings data.
# ------------------------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'159804"
$newSheet.Range("C2").Value = "国寿安保国证创业板中盘精选88ETF"
$newSheet.Range("D2").Value = "'1.10"
$newSheet.Range("E2").Value = "'98.91"
$newSheet.Range("F2").Value = "'1.78"
$newSheet.Range("G2").Value = "'0.0196"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'002872"
$newSheet.Range("C3").Value = "华夏智胜价值成长股票C"
$newSheet.Range("D3").Value = "'2.13"
$newSheet.Range("E3").Value = "'93.39"
$newSheet.Range("F3").Value = "'0.89"
$newSheet.Range("G3").Value = "'0.0190"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'002871"
$newSheet.Range("C4").Value = "华夏智胜价值成长股票A"
$newSheet.Range("D4").Value = "'0.86"
$newSheet.Range("E4").Value = "'93.39"
$newSheet.Range("F4").Value = "'0.89"
$newSheet.Range("G4").Value = "'0.0077"
$newSheet.Range("H4").Value = 7

# ------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: add the 2022-Q3 row at the
#    top of the data and push the older quarters down. Row 5 is new
#    (beyond the original A1:D4 used range) so it has no formatting
#    yet; clone it from row 4 before writing its values so it picks
#    up the same cell style ("A" column bold/border style).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.04
$total.Range("A4:D4").Copy()
$total.Range("A5:D5").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.05

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.2

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.53

# ------------------------------------------------------------------
# 4. Restore the originally-active sheet ("2021-Q4") so the new copy
#    doesn't steal the "tabSelected" flag.
# ------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q4").Activate()
